# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
